$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 46474.297
$ws.Range("J17").Value = 46474.297
$ws.Range("L17").Value = 139422.891
$ws.Range("N17").Value = -139758.891

$ws.Range("H58").Value = 54132
$ws.Range("I58").Value = 287.55554
$ws.Range("J58").Value = 102592
$ws.Range("K58").Value = 862.66662
$ws.Range("L58").Value = 307776
$ws.Range("M58").Value = -712.66662
$ws.Range("N58").Value = -308076

$ws.Range("H82").Value = 3330
$ws.Range("I82").Value = 216.66667
$ws.Range("K82").Value = 650.00001
$ws.Range("M82").Value = -244.00001

$ws.Range("H85").Value = 3330
$ws.Range("I85").Value = 216.66667
$ws.Range("K85").Value = 650.00001
$ws.Range("M85").Value = 753.99999

$ws.Range("H88").Value = 6077.68
$ws.Range("J88").Value = 6098.7144
$ws.Range("L88").Value = 6098.7144
$ws.Range("N88").Value = -6910.7144

$ws.Range("H91").Value = 6077.68
$ws.Range("J91").Value = 6098.7144
$ws.Range("L91").Value = 6098.7144
$ws.Range("N91").Value = -8906.714400000001

$ws.Range("H115").Value = 1095
$ws.Range("I115").Value = 714
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 2142
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -575
$ws.Range("N115").Value = -12134

$ws.Range("H137").Value = 3853.476
$ws.Range("I137").Value = 4362
$ws.Range("J137").Value = 2719.077
$ws.Range("K137").Value = 13086
$ws.Range("L137").Value = 8157.231000000001
$ws.Range("M137").Value = -10536
$ws.Range("N137").Value = -13257.231

$ws.Range("H141").Value = 759376.9399999999
$ws.Range("I141").Value = 6451.7827
$ws.Range("J141").Value = 2333675
$ws.Range("K141").Value = 19355.3481
$ws.Range("L141").Value = 7001025
$ws.Range("M141").Value = -14175.3481
$ws.Range("N141").Value = -7011385

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1535.6923
$ws.Range("I74").Value = 1269.4546
$ws.Range("K74").Value = 1269.4546
$ws.Range("M74").Value = -395.4546

$ws.Range("H77").Value = 1535.6923
$ws.Range("I77").Value = 1269.4546
$ws.Range("K77").Value = 6347.273
$ws.Range("M77").Value = -1979.273

$ws.Range("H104").Value = 28225
$ws.Range("J104").Value = 28225
$ws.Range("L104").Value = 28225
$ws.Range("N104").Value = -35213

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1589.8379
$ws.Range("I20").Value = 1540.3334
$ws.Range("J20").Value = 1654.8125
$ws.Range("K20").Value = 1540.3334
$ws.Range("L20").Value = 1654.8125
$ws.Range("M20").Value = -1293.3334
$ws.Range("N20").Value = -2148.8125

$ws.Range("H107").Value = 2077.8413
$ws.Range("I107").Value = 1644.1702
$ws.Range("J107").Value = 3351.75
$ws.Range("K107").Value = 1644.1702
$ws.Range("L107").Value = 3351.75
$ws.Range("M107").Value = 275.8298
$ws.Range("N107").Value = -7191.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3240.0488
$ws.Range("I31").Value = 1971.4333
$ws.Range("J31").Value = 6699.909
$ws.Range("K31").Value = 1971.4333
$ws.Range("L31").Value = 6699.909
$ws.Range("M31").Value = -1676.4333
$ws.Range("N31").Value = -7289.909

$ws.Range("H34").Value = 3240.0488
$ws.Range("I34").Value = 1971.4333
$ws.Range("J34").Value = 6699.909
$ws.Range("K34").Value = 1971.4333
$ws.Range("L34").Value = 6699.909
$ws.Range("M34").Value = -1769.4333
$ws.Range("N34").Value = -7103.909

$ws.Range("H58").Value = 7044684.5
$ws.Range("I58").Value = 1523.6538
$ws.Range("J58").Value = 26320704
$ws.Range("K58").Value = 1523.6538
$ws.Range("L58").Value = 26320704
$ws.Range("M58").Value = -1320.6538
$ws.Range("N58").Value = -26321110

$ws.Range("H99").Value = 1872.6666
$ws.Range("I99").Value = 1263.4286
$ws.Range("J99").Value = 4005
$ws.Range("K99").Value = 1263.4286
$ws.Range("L99").Value = 4005
$ws.Range("M99").Value = 234.5714
$ws.Range("N99").Value = -7001

$ws.Range("H126").Value = 1872.6666
$ws.Range("I126").Value = 1263.4286
$ws.Range("J126").Value = 4005
$ws.Range("K126").Value = 3790.2858
$ws.Range("L126").Value = 12015
$ws.Range("M126").Value = -1320.2858
$ws.Range("N126").Value = -16955

$ws.Range("H132").Value = 2093.4443
$ws.Range("I132").Value = 1597.8529
$ws.Range("K132").Value = 4793.5587
$ws.Range("M132").Value = -2263.5587

$ws.Range("H134").Value = 1999.0605
$ws.Range("I134").Value = 1113.5
$ws.Range("K134").Value = 3340.5
$ws.Range("M134").Value = -805.5

$ws.Range("H136").Value = 7044684.5
$ws.Range("I136").Value = 1523.6538
$ws.Range("J136").Value = 26320704
$ws.Range("K136").Value = 4570.9614
$ws.Range("L136").Value = 78962112
$ws.Range("M136").Value = -2020.9614
$ws.Range("N136").Value = -78967212

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 7748.3335
$ws.Range("I87").Value = 4564.1665
$ws.Range("J87").Value = 14116.667
$ws.Range("K87").Value = 13692.4995
$ws.Range("L87").Value = 42350.001
$ws.Range("M87").Value = -12444.4995
$ws.Range("N87").Value = -44846.001

$ws.Range("H90").Value = 7748.3335
$ws.Range("I90").Value = 4564.1665
$ws.Range("J90").Value = 14116.667
$ws.Range("K90").Value = 41077.4985
$ws.Range("L90").Value = 127050.003
$ws.Range("M90").Value = -34837.4985
$ws.Range("N90").Value = -139530.003

$ws.Range("H120").Value = 17420.818
$ws.Range("I120").Value = 13257.5
$ws.Range("K120").Value = 39772.5
$ws.Range("M120").Value = -34934.5

$ws.Range("H140").Value = 16669533
$ws.Range("I140").Value = 33334306
$ws.Range("K140").Value = 100002918
$ws.Range("M140").Value = -99997738

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 4999
$ws.Range("I52").Value = 4999
$ws.Range("K52").Value = 4999
$ws.Range("M52").Value = -4740

$ws.Range("H97").Value = 946.04
$ws.Range("I97").Value = 431.57895
$ws.Range("J97").Value = 2575.1667
$ws.Range("K97").Value = 431.57895
$ws.Range("L97").Value = 2575.1667
$ws.Range("M97").Value = 64.42104999999998
$ws.Range("N97").Value = -3567.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1314.4814
$ws.Range("I46").Value = 533.3333
$ws.Range("J46").Value = 1412.125
$ws.Range("K46").Value = 533.3333
$ws.Range("L46").Value = 1412.125
$ws.Range("M46").Value = -345.3333
$ws.Range("N46").Value = -1788.125

$ws.Range("H93").Value = 1325.5625
$ws.Range("I93").Value = 729.2143
$ws.Range("K93").Value = 729.2143
$ws.Range("M93").Value = 518.7857

$ws.Range("H132").Value = 1581.48
$ws.Range("I132").Value = 992.8982999999999
$ws.Range("J132").Value = 3751.875
$ws.Range("K132").Value = 2978.6949
$ws.Range("L132").Value = 11255.625
$ws.Range("M132").Value = -448.6949
$ws.Range("N132").Value = -16315.625

$ws.Range("H136").Value = 1585.7059
$ws.Range("I136").Value = 1105.8914
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 3317.6742
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -767.6741999999999
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8465.831
$ws.Range("I132").Value = 1960.0625
$ws.Range("J132").Value = 19234
$ws.Range("K132").Value = 5880.1875
$ws.Range("L132").Value = 57702
$ws.Range("M132").Value = -3350.1875
$ws.Range("N132").Value = -62762

$ws.Range("H133").Value = 39563.332
$ws.Range("J133").Value = 39563.332
$ws.Range("L133").Value = 39563.332
$ws.Range("N133").Value = -49683.332

$ws.Range("H136").Value = 1428.7354
$ws.Range("I136").Value = 831.6875
$ws.Range("J136").Value = 1959.4445
$ws.Range("K136").Value = 2495.0625
$ws.Range("L136").Value = 5878.333500000001
$ws.Range("M136").Value = 54.9375
$ws.Range("N136").Value = -10978.3335
